# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert "Georgia" as a new row right before "Siria" (row 131), pushing
#    "Siria" and "Lituania" down by one row, and remove the old "Georgia"
#    row that is now a duplicate further down (originally row 133, now 134
#    after the insert).
# ---------------------------------------------------------------------------
$ws.Rows.Item(131).Insert()

$ws.Cells.Item(131, 1).Value = "Georgia"
$ws.Cells.Item(131, 2).Value = 3913
$ws.Cells.Item(131, 3).Value = 218
$ws.Cells.Item(131, 4).Value = 1574
$ws.Cells.Item(131, 5).Value = 2317
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 2
$ws.Cells.Item(131, 8).Value = 22

# Remove the now-duplicate old "Georgia" row (pushed down to row 134).
$ws.Rows.Item(134).Delete()

# ---------------------------------------------------------------------------
# 2) Swap "Islas Malvinas" (row 214) and "Montserrat" (row 215) so that
#    "Montserrat" now comes first.
# ---------------------------------------------------------------------------
$ws.Cells.Item(214, 1).Value = "Montserrat"
$ws.Cells.Item(214, 2).Value = 13
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 1

$ws.Cells.Item(215, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 2).Value = 13
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 5).Value = 0
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 0

# ---------------------------------------------------------------------------
# 3) Refresh case counts for a handful of countries.
# ---------------------------------------------------------------------------
# Israel (row 27)
$ws.Cells.Item(27, 2).Value = 192579
$ws.Cells.Item(27, 3).Value = 1650
$ws.Cells.Item(27, 4).Value = 140743
$ws.Cells.Item(27, 5).Value = 50563

# Barein (row 54)
$ws.Cells.Item(54, 5).Value = 6901
$ws.Cells.Item(54, 7).Value = 1
$ws.Cells.Item(54, 8).Value = 225

# Suiza (row 61)
$ws.Cells.Item(61, 4).Value = 41800
$ws.Cells.Item(61, 5).Value = 6528

# Armenia (row 63)
$ws.Cells.Item(63, 2).Value = 47667
$ws.Cells.Item(63, 3).Value = 115
$ws.Cells.Item(63, 4).Value = 42676
$ws.Cells.Item(63, 5).Value = 4053
$ws.Cells.Item(63, 7).Value = 2
$ws.Cells.Item(63, 8).Value = 938

# ---------------------------------------------------------------------------
# 4) Update the "last updated" timestamp banner.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 22 de Septiembre de 2020 a las 09:11"
